$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2666666666666667
$ws.Range("C2").Value = 0.4
$ws.Range("P2").Value = 0.2
$ws.Range("S2").Value = 0.1333333333333333
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.5
$ws.Range("P4").Value = 0.8
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.0625
$ws.Range("F6").Value = 0.125
$ws.Range("J6").Value = 0.25
$ws.Range("Q6").Value = 0.125
$ws.Range("S6").Value = 0.4375
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("Q7").Value = 0.05555555555555555
$ws.Range("R7").Value = 0.1111111111111111
$ws.Range("B8").Value = 0.03703703703703703
$ws.Range("F8").Value = 0.03703703703703703
$ws.Range("J8").Value = 0.1111111111111111
$ws.Range("Q8").Value = 0.1111111111111111
$ws.Range("R8").Value = 0.03703703703703703
$ws.Range("S8").Value = 0.6666666666666666
$ws.Range("B9").Value = 0.03846153846153846
$ws.Range("D9").Value = 0.03846153846153846
$ws.Range("F9").Value = 0.03846153846153846
$ws.Range("J9").Value = 0.1153846153846154
$ws.Range("Q9").Value = 0.1923076923076923
$ws.Range("R9").Value = 0.03846153846153846
$ws.Range("S9").Value = 0.5384615384615384
$ws.Range("B10").Value = 0.09523809523809523
$ws.Range("D10").Value = 0.04761904761904762
$ws.Range("F10").Value = 0.09523809523809523
$ws.Range("J10").Value = 0.1428571428571428
$ws.Range("Q10").Value = 0.1666666666666667
$ws.Range("R10").Value = 0.07142857142857142
$ws.Range("S10").Value = 0.3809523809523809
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.03333333333333333
$ws.Range("K11").Value = 0.2333333333333333
$ws.Range("L11").Value = 0.5
$ws.Range("S11").Value = 0.06666666666666667
$ws.Range("G12").Value = 0.7333333333333333
$ws.Range("J12").Value = 0.06666666666666667
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.1
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.3
$ws.Range("J15").Value = 0.4
$ws.Range("H16").Value = 0.3333333333333333
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("O16").Value = 0.1111111111111111
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("H17").Value = 0.2692307692307692
$ws.Range("I17").Value = 0.1538461538461539
$ws.Range("J17").Value = 0.2692307692307692
$ws.Range("K17").Value = 0.1153846153846154
$ws.Range("S17").Value = 0.1923076923076923
$ws.Range("H18").Value = 0.2222222222222222
$ws.Range("I18").Value = 0.3333333333333333
$ws.Range("J18").Value = 0.2222222222222222
$ws.Range("S18").Value = 0.2222222222222222
$ws.Range("F19").Value = 0.008
$ws.Range("H19").Value = 0.112
$ws.Range("I19").Value = 0.12
$ws.Range("J19").Value = 0.328
$ws.Range("K19").Value = 0.152
$ws.Range("M19").Value = 0.024
$ws.Range("N19").Value = 0.008
$ws.Range("O19").Value = 0.056
$ws.Range("S19").Value = 0.192
